$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("J2").Value = 3.4
$ws.Range("S2").Value = 1.62
$ws.Range("T2").Value = 2.2
$ws.Range("U2").Value = 2.25
$ws.Range("V2").Value = 1.57
$ws.Range("W2").Value = 6
$ws.Range("Z2").Value = 23
$ws.Range("AE2").Value = 19
$ws.Range("AO2").Value = 15
$ws.Range("AT2").Value = 2.2
$ws.Range("AU2").Value = 9.5
$ws.Range("AY2").Value = 34
$ws.Range("BB2").Value = 500

# Row 3 updates
$ws.Range("Q3").Value = 1.67
$ws.Range("R3").Value = 2.15

# Row 4 updates
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 3.2
$ws.Range("I4").Value = 3.9
$ws.Range("J4").Value = 2.75
$ws.Range("L4").Value = 4.33
$ws.Range("Q4").Value = 2.08
$ws.Range("R4").Value = 1.73
$ws.Range("AD4").Value = 6
$ws.Range("AI4").Value = 19
$ws.Range("AJ4").Value = 13
$ws.Range("AY4").Value = 29
$ws.Range("AZ4").Value = 67
